# Update to published CDA FHIR logical model with patches #241
#
# Sheet "Metadata": bump Version / Date / Contact values.
# Sheet "Include from EntityNameUse": insert a new "ASGN" concept row
# between "SNDX" and "ABC".

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$wsMeta = $wb.Worksheets.Item(1)

$wsMeta.Cells.Item(3, 2).Value = "2.0.0-sd-202406-matchbox-patch"
$wsMeta.Cells.Item(8, 2).Value = "2024-06-19T17:47:42+02:00"
$wsMeta.Cells.Item(10, 2).Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- Concept sheet ----------------------------------------------------
$wsConcept = $wb.Worksheets.Item(2)

# Insert a new row 11 (pushes ABC/SYL/IDE/blank/System URI rows down by one)
$wsConcept.Rows.Item(11).Insert()

# Match the formatting of the surrounding concept rows (copy from row 10)
$wsConcept.Range("A10:B10").Copy()
$wsConcept.Range("A11:B11").PasteSpecial(-4122)

# New concept code
$wsConcept.Cells.Item(11, 1).Value = "ASGN"

# Row 14 now holds "IDE" (shifted down from row 13); its B cell has always
# been empty, but the copy/paste above leaves a stray empty-string value in
# it - clear it back out to an unset cell, same as all the other code rows.
$wsConcept.Cells.Item(14, 2).ClearContents()
